$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: E6 becomes a formula (60*30) instead of the literal 90,
# and a new J6 cell is added with the "probability of total" over 60 trials.
$ws.Range("E6").Formula = "=60*30"
$ws.Range("J6").Formula = "=1-(1-C6)^60"

# Row 7: MTTH event (Interval = 35*12, count = 1)
$ws.Range("E7").Formula = "=35*12"
$ws.Range("F7").Value = 1
$ws.Range("C7").Formula = "=1-(0.5)^(F7/E7)"
$ws.Range("D7").Formula = "=F7*(1-0.5^(1/E7))"

# Row 8: MTTH event (Interval = 600, count = 1)
$ws.Range("E8").Value = 600
$ws.Range("F8").Value = 1
$ws.Range("C8").Formula = "=1-(0.5)^(F8/E8)"
$ws.Range("D8").Formula = "=F8*(1-0.5^(1/E8))"

# Row 9: MTTH event (Interval = 150, count = 10)
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 10
$ws.Range("C9").Formula = "=1-(0.5)^(F9/E9)"
$ws.Range("D9").Formula = "=F9*(1-0.5^(1/E9))"

# Row 10: MTTH event (Interval = 180, count = 10)
$ws.Range("E10").Value = 180
$ws.Range("F10").Value = 10
$ws.Range("C10").Formula = "=1-(0.5)^(F10/E10)"
$ws.Range("D10").Formula = "=F10*(1-0.5^(1/E10))"

# Update selection to match the target state
$ws.Range("D15").Select()
